$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new tracker entry as row 49 (A1:F48 -> A1:F49)
$ws.Range("A49").Value = "G2"
$ws.Range("B49").Value = "Workout"

# Date cell: match the existing date-formatted column (style used by C2:C48)
$ws.Range("C49").Value = 45907
$ws.Range("C49").NumberFormat = $ws.Range("C48").NumberFormat

$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
